# Update the per-seed metric results (columns B:K, rows 2-16) with the
# newly-recomputed / rounded values, as part of introducing a visualization
# of results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @(0.75, 0.69, 0.73, 0.93, 0.99, 0.87, 0.9, 0.93, 0.87, 0.85)
    3 = @(0.76, 0.68, 0.73, 0.95, 0.99, 0.89, 0.9, 0.95, 0.88, 0.86)
    4 = @(0.8, 0.64, 0.76, 0.87, 0.99, 0.86, 0.92, 0.95, 0.85, 0.85)
    5 = @(0.75, 0.68, 0.73, 0.89, 0.97, 0.86, 0.9, 0.95, 0.87, 0.84)
    6 = @(0.75, 0.7, 0.72, 0.94, 0.98, 0.87, 0.91, 0.95, 0.84, 0.85)
    7 = @(0.76, 0.67, 0.7, 0.85, 0.98, 0.87, 0.89, 0.95, 0.85, 0.84)
    8 = @(0.77, 0.7, 0.71, 0.9, 1, 0.89, 0.9, 0.96, 0.85, 0.85)
    9 = @(0.75, 0.67, 0.69, 0.87, 0.99, 0.87, 0.9, 0.96, 0.86, 0.84)
    10 = @(0.74, 0.64, 0.69, 0.83, 0.98, 0.89, 0.9, 0.95, 0.85, 0.83)
    11 = @(0.77, 0.69, 0.69, 0.93, 0.97, 0.87, 0.92, 0.96, 0.88, 0.85)
    12 = @(0.72, 0.68, 0.71, 0.97, 0.97, 0.86, 0.9, 0.95, 0.82, 0.84)
    13 = @(0.77, 0.68, 0.72, 0.93, 0.97, 0.9, 0.91, 0.96, 0.88, 0.86)
    14 = @(0.75, 0.65, 0.68, 0.94, 0.98, 0.87, 0.91, 0.96, 0.86, 0.84)
    15 = @(0.75, 0.66, 0.69, 0.96, 0.98, 0.88, 0.91, 0.95, 0.87, 0.85)
    16 = @(0.76, 0.67, 0.71, 0.91, 0.98, 0.88, 0.91, 0.95, 0.86, 0.85)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $col = 2
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
